# timelog * user interaction box * working on rubics cube scene

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 12 with a new timelog entry
$ws.Cells.Item(12, 1).Value = 41790        # A12 - date
$ws.Cells.Item(12, 2).Value = 0.125        # B12 - from
$ws.Cells.Item(12, 3).Value = 0.20833333333333334  # C12 - to
$ws.Cells.Item(12, 5).Value = "physics scene refactoring, rubics cube scene"  # E12 - activity

# Update the selection to reflect the active cell after editing
$ws.Range("E12").Select()
